$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new data set introduces two additional weekly price records for
# "Brocoli" at "Vega Monumental Concepcion" (fecha = 45106). They are
# inserted right before the existing row that used to be row 492, which
# pushes all the following rows (492..527) down by two, ending at 529.
$ws.Rows("492:493").Insert()

# --- New row 492 (Primera) ---
$ws.Range("A492").Value = 11
$ws.Range("B492").Value = "Vega Monumental Concepción"
$ws.Range("C492").Value = "Bíobío"
$ws.Range("D492").Value = 45106
$ws.Range("D492").NumberFormat = $ws.Range("D491").NumberFormat()
$ws.Range("E492").Value = 8
$ws.Range("F492").Value = 100112023
$ws.Range("G492").Value = "Brócoli"
$ws.Range("H492").Value = "Sin especificar"
$ws.Range("I492").Value = "Primera"
$ws.Range("J492").Value = 2000
$ws.Range("K492").Value = 700
$ws.Range("L492").Value = 800
$ws.Range("M492").Value = 750
$ws.Range("N492").Value = "$/unidad"
$ws.Range("O492").Value = "Región Metropolitana"
$ws.Range("P492").Value = 750
$ws.Range("Q492").Value = 1
$ws.Range("R492").Value = "Hortaliza"

# --- New row 493 (Segunda) ---
$ws.Range("A493").Value = 11
$ws.Range("B493").Value = "Vega Monumental Concepción"
$ws.Range("C493").Value = "Bíobío"
$ws.Range("D493").Value = 45106
$ws.Range("D493").NumberFormat = $ws.Range("D491").NumberFormat()
$ws.Range("E493").Value = 8
$ws.Range("F493").Value = 100112023
$ws.Range("G493").Value = "Brócoli"
$ws.Range("H493").Value = "Sin especificar"
$ws.Range("I493").Value = "Segunda"
$ws.Range("J493").Value = 1000
$ws.Range("K493").Value = 600
$ws.Range("L493").Value = 600
$ws.Range("M493").Value = 600
$ws.Range("N493").Value = "$/unidad"
$ws.Range("O493").Value = "Región Metropolitana"
$ws.Range("P493").Value = 600
$ws.Range("Q493").Value = 1
$ws.Range("R493").Value = "Hortaliza"
